# Auto-generated edit script: updates crypto price/volume table cells
# per the commit 'Updated cryptos list on Wed Oct  4 05:50:31 UTC 2023 with GitHub Actions'

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.394.05'
$ws.Range("E2").Value = '  -0.72%  '
$ws.Range("D3").Value = '1.639.49'
$ws.Range("E3").Value = '  -1.71%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.529'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.62%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.06'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.63%  '
$ws.Range("E9").Value = '  -2.20%  '
$ws.Range("E10").Value = '  -2.11%  '
$ws.Range("E11").Value = '  +1.08%  '
$ws.Range("D12").Value = '1.870.70'
$ws.Range("E12").Value = '  -1.74%  '
$ws.Range("D13").Value = '1.651.93'
$ws.Range("E13").Value = '  -1.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.559'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.32'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.96%  '
$ws.Range("D17").Value = '27.366.62'
$ws.Range("E17").Value = '  -0.84%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.98'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.60%  '
$ws.Range("E19").Value = '  -1.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("E22").Value = '  -3.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.53'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.68%  '
$ws.Range("E24").Value = '  -0.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.51'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.45%  '
$ws.Range("E26").Value = '  -3.30%  '
$ws.Range("E27").Value = '  +1.46%  '
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.55'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.52%  '
$ws.Range("E30").Value = '  -4.05%  '
$ws.Range("E31").Value = '  -3.43%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.28'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.16%  '
$ws.Range("E33").Value = '  +0.09%  '
$ws.Range("D34").Value = '1.408.87'
$ws.Range("E34").Value = '  -3.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.58'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.42%  '
$ws.Range("E36").Value = '  -0.35%  '
$ws.Range("E37").Value = '  -1.93%  '
$ws.Range("E38").Value = '  -5.40%  '
$ws.Range("E39").Value = '  -3.61%  '
$ws.Range("E40").Value = '  +1.21%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.53'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.08%  '
$ws.Range("B43").Value = 'mCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.46'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.00%  '
$ws.Range("E44").Value = '  +0.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.793'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.60%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.47'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.36%  '
$ws.Range("D47").Value = '1.780.75'
$ws.Range("E47").Value = '  -1.63%  '
$ws.Range("E48").Value = '  -4.64%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '87.36'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.27%  '
$ws.Range("E50").Value = '  -2.67%  '
$ws.Range("E51").Value = '  -3.69%  '
